# CYBER360-8.3-Inheritance.pptx — "fix typos, ch 5, 6, 8"
#
# The only substantive content edit is on the title slide (Slide 1):
# the Title placeholder text "CIT 361/CYBER 360: Advanced Scripting"
# is corrected to "CYBER 360: Advanced Scripting".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$titleShape = $s.Shapes.Item("Title 1")
$titleShape.TextFrame.TextRange.Text = "CYBER 360: Advanced Scripting"
